$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the NPV values for each scenario (B2:B4)
$ws.Range("B2").Value = "$1197 Billion"
$ws.Range("B3").Value = "$1023 Billion"
$ws.Range("B4").Value = "$994 Billion"

# Move selection to B5, matching the updated sheet view state
$ws.Range("B5").Select()
